$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark from its old location
#    (it currently sits between the "a" run and the "." run in the
#    "suitability" paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. In the dictionary paragraph, replace the long-winded description of how
#    words are stored/linked with the shorter "têm uma" phrasing.
$d.Content.Find.Execute(
    "devem ser armazenadas e associadas a um link que redireciona para uma página que conterá a explicação",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "têm uma explicação", 2) | Out-Null

# 3. Split the run containing "têm uma" off from its neighbours (so the new
#    bookmark can sit in its own run boundary) by toggling a character
#    property on and back off again, which forces Word to break the run
#    without actually changing any visible formatting.
$find2 = $d.Content.Find
$find2.Execute("têm uma") | Out-Null
$temEspecie = $find2.Parent
$temEspecie.Bold = 1
$temEspecie.Bold = 0

# 4. Re-insert the "_GoBack" bookmark as a collapsed range immediately
#    before "têm uma".
$bookmarkRange = $d.Range($temEspecie.Start, $temEspecie.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# 5. In the investment-simulation paragraph, the text run was previously
#    split in two around a <w:lastRenderedPageBreak/> marker. Re-running a
#    find/replace over that same text causes Word to rebuild the run and
#    drop the stale page-break marker, merging the text back into one run.
$d.Content.Find.Execute(
    "tipo de investimento que o usuário escolheu",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "tipo de investimento que o usuário escolheu", 2) | Out-Null
